$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Blad1")

# Week 15 table (rows 113-117), Maandag row (113): 8 -> 4
$ws.Range("B113:C113").Value = 4
$ws.Range("E113:G113").Value = 4

# Week 15 table, Donderdag row (116): 4 -> 2
$ws.Range("B116:C116").Value = 2
$ws.Range("E116:G116").Value = 2

# Week 15 table, Vrijdag row (117) was empty, now filled in
$ws.Range("B117").Value = 4
$ws.Range("C117").Value = 4
$ws.Range("D117").Value = 0
$ws.Range("E117").Value = 4
$ws.Range("F117").Value = 4
$ws.Range("G117").Value = 4

# Highlight like the other filled-in days: blue for present, red for the 0-hour day
$ws.Range("C117").Interior.Color = 12611584
$ws.Range("D117").Interior.Color = 255
$ws.Range("E117").Interior.Color = 12611584
$ws.Range("F117").Interior.Color = 12611584
$ws.Range("G117").Interior.Color = 12611584

# Update the selected cell to reflect the saved view state
$ws.Range("K16").Select()
